$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19: fermentation improvement value change
$ws.Range("E19").Value = 0.73

# Row 20: fermentation improvement value change
$ws.Range("E20").Value = 68

# Row 24: change value and rewrite formulas to reference E24 instead of hardcoded 0.05
$ws.Range("E24").Value = 0.048
$ws.Range("G24").Formula = "=E24*0.0463/0.2087"
$ws.Range("I24").Formula = "=E24*0.34/0.2087"

# Update the selection to H27 as in the diff
$ws.Range("H27").Select()
